# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

$ws.Range("D2").Value = "36.705.93"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "2.078.20"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue $ws.Range("D5") "0.715"
$ws.Range("E5").Value = "  +7.35%  "
Set-TextValue $ws.Range("D6") "244.12"
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("E7").Value = "  +0.06%  "
Set-TextValue $ws.Range("D8") "52.43"
$ws.Range("E8").Value = "  -7.99%  "
Set-TextValue $ws.Range("D9") "58.77"
$ws.Range("E9").Value = "  -2.37%  "
Set-TextValue $ws.Range("D10") "0.364"
$ws.Range("E10").Value = "  -5.52%  "
Set-TextValue $ws.Range("D11") "0.0756"
$ws.Range("E11").Value = "  -3.75%  "
Set-TextValue $ws.Range("D12") "0.109"
$ws.Range("E12").Value = "  +0.66%  "
Set-TextValue $ws.Range("D13") "0.900"
$ws.Range("E13").Value = "  -1.45%  "
Set-TextValue $ws.Range("D14") "14.59"
$ws.Range("E14").Value = "  -10.40%  "
$ws.Range("D15").Value = "2.386.80"
$ws.Range("E15").Value = "  +1.43%  "
Set-TextValue $ws.Range("D16") "5.43"
$ws.Range("E16").Value = "  -6.14%  "
$ws.Range("D17").Value = "2.144.04"
$ws.Range("E17").Value = "  +4.28%  "
$ws.Range("D18").Value = "36.663.89"
$ws.Range("E18").Value = "  -1.46%  "
Set-TextValue $ws.Range("D19") "16.70"
$ws.Range("E19").Value = "  -10.53%  "
Set-TextValue $ws.Range("D20") "72.44"
$ws.Range("E20").Value = "  -3.14%  "
$ws.Range("D21").Value = "0.0₃0868"
$ws.Range("E21").Value = "  -3.76%  "
Set-TextValue $ws.Range("D22") "5.37"
$ws.Range("E22").Value = "  -2.31%  "
Set-TextValue $ws.Range("D23") "237.71"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  -0.09%  "
Set-TextValue $ws.Range("D25") "2.38"
$ws.Range("E25").Value = "  -4.43%  "
Set-TextValue $ws.Range("D26") "9.52"
$ws.Range("E26").Value = "  -1.67%  "
Set-TextValue $ws.Range("D27") "2.14"
$ws.Range("E27").Value = "  -2.04%  "
Set-TextValue $ws.Range("D28") "165.29"
$ws.Range("E28").Value = "  -3.02%  "
Set-TextValue $ws.Range("D29") "20.49"
$ws.Range("E29").Value = "  +1.18%  "
Set-TextValue $ws.Range("D30") "0.131"
$ws.Range("E30").Value = "  +4.27%  "
Set-TextValue $ws.Range("D31") "5.15"
$ws.Range("E31").Value = "  -1.38%  "
Set-TextValue $ws.Range("D32") "1.15"
$ws.Range("E32").Value = "  -2.86%  "
Set-TextValue $ws.Range("D33") "4.64"
$ws.Range("E33").Value = "  +0.76%  "
Set-TextValue $ws.Range("D34") "0.0600"
$ws.Range("E34").Value = "  -4.11%  "
Set-TextValue $ws.Range("D35") "2.38"
$ws.Range("E35").Value = "  +4.87%  "
$ws.Range("E36").Value = "  +0.29%  "
Set-TextValue $ws.Range("D37") "1.82"
$ws.Range("E37").Value = "  +2.41%  "
Set-TextValue $ws.Range("D38") "0.0816"
$ws.Range("E38").Value = "  -7.84%  "
Set-TextValue $ws.Range("D39") "1.26"
$ws.Range("E39").Value = "  -6.90%  "
$ws.Range("B40").Value = "THORChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D40") "4.85"
$ws.Range("E40").Value = "  -8.16%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D41") "0.0218"
$ws.Range("E41").Value = "  -2.63%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D42") "1.14"
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D43") "0.0947"
$ws.Range("E43").Value = "  -4.99%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D44") "2.87"
$ws.Range("E44").Value = "  -6.86%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D45") "95.02"
$ws.Range("E45").Value = "  -1.71%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.390.88"
$ws.Range("E46").Value = "  +9.06%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D47") "7.55"
$ws.Range("E47").Value = "  +9.96%  "
Set-TextValue $ws.Range("D48") "15.69"
$ws.Range("E48").Value = "  -11.65%  "
Set-TextValue $ws.Range("D49") "2.39"
$ws.Range("E49").Value = "  -3.03%  "
Set-TextValue $ws.Range("D50") "2.87"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("D51").Value = "2.273.85"
$ws.Range("E51").Value = "  +1.52%  "
